$d = $word.ActiveDocument

$replacements = @(
    @("175×6=", "718×2="),
    @("749×2=", "432×5="),
    @("315×3=", "567×7="),
    @("477×6=", "915×7="),
    @("423×2=", "265×7="),
    @("890×8=", "875×7="),
    @("956×4=", "213×6="),
    @("161×2=", "388×9="),
    @("373×5=", "649×6="),
    @("736×3=", "403×9="),
    @("860×5=", "215×7="),
    @("327×7=", "674×6="),
    @("938×4=", "461×8="),
    @("713×3=", "572×7="),
    @("451×6=", "409×9="),
    @("811×4=", "325×5="),
    @("707×2=", "364×7="),
    @("569×9=", "902×3="),
    @("206×6=", "471×6="),
    @("867×7=", "852×8="),
    @("151×5=", "309×4="),
    @("324×8=", "880×8="),
    @("765×4=", "296×7="),
    @("733×6=", "665×9="),
    @("623×9=", "440×3=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
